# Generate Report for Handoff
# Replace the old handoff-file UUID/basename ("3eaee533-2d7b-47ac-b4de-1074a8cbed8b")
# with a new one ("cf39f490-87b4-4cea-9542-1190327a7289"), refresh the dependent
# handoff-file names (new content hash "2d41a07744acab15c4e3000818a0d0244b200ba0")
# and bump the "Latest Handoff Datetime" stamps, across all three sheets
# (Overview, zh-cn, de-de). Hyperlink targets (URLs) are left pointing at the
# same places as before - only the cell text / display text changes.

$wb = $excel.ActiveWorkbook

$oldMdName  = "3eaee533-2d7b-47ac-b4de-1074a8cbed8b.md"
$newMdName  = "cf39f490-87b4-4cea-9542-1190327a7289.md"

$oldZhXlf = "3eaee533-2d7b-47ac-b4de-1074a8cbed8b.91feb01bb7403268261094acbde4a13f80c1942d.zh-cn.xlf"
$newZhXlf = "cf39f490-87b4-4cea-9542-1190327a7289.2d41a07744acab15c4e3000818a0d0244b200ba0.zh-cn.xlf"

$oldDeXlf = "3eaee533-2d7b-47ac-b4de-1074a8cbed8b.91feb01bb7403268261094acbde4a13f80c1942d.de-de.xlf"
$newDeXlf = "cf39f490-87b4-4cea-9542-1190327a7289.2d41a07744acab15c4e3000818a0d0244b200ba0.de-de.xlf"

$oldZhDate = "2016-02-22 05:05:32"
$newZhDate = "2016-02-22 05:06:34"

$oldDeDate = "2016-02-22 05:05:47"
$newDeDate = "2016-02-22 05:06:48"

$mdTarget = "https://github.com/OpenLocalizationTest/oltest/blob/39b8b21f39222b28f99365acffdb706f89fe02dd/e2e/$oldMdName"
$configTarget = "https://github.com/OpenLocalizationTest/oltest/blob/39b8b21f39222b28f99365acffdb706f89fe02dd/.localization-config"
$zhXlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0df4bddd0cb96de4af5b7b3c0f15cb450c037865/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$oldZhXlf"
$deXlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d26acbb191882093eb410b7b45fc4229388dffca/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$oldDeXlf"

# ---- Sheet "Overview": A2 = handoff markdown file name/link ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdTarget, "", "", $newMdName) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $configTarget, "", "", ".localization-config") | Out-Null

# ---- Sheet "zh-cn": A2 = handoff markdown; C2 = latest handoff xlf; D2 = datetime ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newMdName
$wsZh.Range("C2").Value = $newZhXlf
$wsZh.Range("D2").Value = $newZhDate
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $mdTarget, "", "", $newMdName) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), $zhXlfTarget, "", "", $newZhXlf) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $configTarget, "", "", ".localization-config") | Out-Null

# ---- Sheet "de-de": A2 = handoff markdown; C2 = latest handoff xlf; D2 = datetime ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newMdName
$wsDe.Range("C2").Value = $newDeXlf
$wsDe.Range("D2").Value = $newDeDate
$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $mdTarget, "", "", $newMdName) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), $deXlfTarget, "", "", $newDeXlf) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $configTarget, "", "", ".localization-config") | Out-Null
